# fix: changes from review
# Adds a new "AnniversarySheet" worksheet (after the existing "BirthdaySheet")
# with a Name 1 / Name 2 / Anniversary table, and tweaks column widths.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- New worksheet, placed right after BirthdaySheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "AnniversarySheet"

# --- Header row ---
$ws2.Range("A1").Value = "Name 1"
$ws2.Range("B1").Value = "Name 2"
$ws2.Range("C1").Value = "Anniversary"

# --- Data row ---
$ws2.Range("A2").Value = "Mickey Mouse"
$ws2.Range("B2").Value = "Mini Mouse"
$ws2.Range("C2").Value = 4019
$ws2.Range("C2").NumberFormat = "mm/dd/yy"

# --- Reuse BirthdaySheet's existing cell styles so no duplicate style
#     entries get minted in styles.xml (copy formats only) ---
$ws1.Range("A1").Copy()
$ws2.Range("A1:A3").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("A1").Copy()
$ws2.Range("B1").PasteSpecial(-4122)      # xlPasteFormats

$ws1.Range("B2").Copy()
$ws2.Range("B2:B3").PasteSpecial(-4122)   # xlPasteFormats

# --- Column widths ---
$ws1.Columns.Item(1).ColumnWidth = 27.68

$ws2.Columns.Item(1).ColumnWidth = 27.82
$ws2.Columns.Item(2).ColumnWidth = 24.49

# --- Selection / active sheet state ---
$ws2.Range("A3").Select()
$ws1.Select()
$ws1.Range("A1").Select()
